$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 227, shifting existing rows 227-237 down to 228-238.
$ws.Rows("227:227").Insert()

# Populate the newly-inserted row 227 with the new weekly record.
$ws.Range("A227").Value = 3
$ws.Range("B227").Value = "Femacal de La Calera"
$ws.Range("C227").Value = "Coquimbo"
$ws.Range("D227").Value = 44509
$ws.Range("D227").NumberFormat = $ws.Range("D228").NumberFormat
$ws.Range("E227").Value = 5
$ws.Range("F227").Value = 100112040
$ws.Range("G227").Value = "Cilantro"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 270
$ws.Range("K227").Value = 2000
$ws.Range("L227").Value = 2300
$ws.Range("M227").Value = 2167
$ws.Range("N227").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O227").Value = "Provincia de Quillota"
$ws.Range("P227").Value = 722
$ws.Range("Q227").Value = 3
$ws.Range("R227").Value = "Hortaliza"
